{"js": "// Clarify language around impact at Avanade: the first bullet under the\n// \"Analyst, Machine Learning and Azure Cloud Enablement\" role changes from\n//   \"...machine learning for a cloud resource optimizer creating yearly\n//    client savings of $600K+\"\n// to\n//   \"...machine learning for an optimizer yielding client savings of\n//    25-40% total Azure spend\"\n\nconst body = context.document.body;\n\n// Locate the exact original sentence fragment (unique in the document) so\n// we don't accidentally touch anything else.\nconst searchResults = body.search(\n  \"for a cloud resource optimizer creating yearly client savings of $600K+\",\n  { matchCase: true }\n);\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the Avanade impact sentence to update.\");\n}\n\n// Replace it in place; formatting (font, size, color, etc.) of the\n// surrounding run carries over to the newly inserted text.\nsearchResults.items[0].insertText(\n  \"for an optimizer yielding client savings of 25-40% total Azure spend\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Clarify language around impact at Avanade: the first bullet under the\n# \"Analyst, Machine Learning and Azure Cloud Enablement\" role changes from\n#   \"...machine learning for a cloud resource optimizer creating yearly\n#    client savings of $600K+\"\n# to\n#   \"...machine learning for an optimizer yielding client savings of\n#    25-40% total Azure spend\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"for a cloud resource optimizer creating yearly client savings of `$600K+\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"for an optimizer yielding client savings of 25-40% total Azure spend\",\n    2\n)\n\nif (-not $found) {\n    throw \"Could not find the Avanade impact sentence to update.\"\n}\n"}
